$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (StreetTRACKS Gold Shares / GLD)
$ws.Range("D2").Value = 387.49
$ws.Range("E2").Value = 56.7
$ws.Range("F2").Value = 1.14
$ws.Range("K2").Value = 67.3
$ws.Range("N2").Value = 53.62998959737769

# Row 3 (Newmont Corporation / NEM)
$ws.Range("D3").Value = 90.59999999999999
$ws.Range("E3").Value = 52.4
$ws.Range("F3").Value = 0.08
$ws.Range("G3").Value = 60
$ws.Range("K3").Value = 66.09999999999999
$ws.Range("N3").Value = 53.62998959737769

# Row 4 (Gold Feb 26 / GC=F)
$ws.Range("D4").Value = 4242
$ws.Range("E4").Value = 71.8
$ws.Range("F4").Value = 4.45
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 63
$ws.Range("K4").Value = 53.3
$ws.Range("N4").Value = 53.62998959737769
